$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New date column header (next day after 31-jul)
$ws.Range("AP1").Value = "01-ago"

# New column values for each data row
$ws.Range("AP2").Value = 0
$ws.Range("AP3").Value = 17.739341569465441
$ws.Range("AP4").Value = 17.25967898842293
$ws.Range("AP5").Value = 23.484576240675192
$ws.Range("AP6").Value = 0
$ws.Range("AP7").Value = 11.361505408100504
$ws.Range("AP8").Value = 10.854212682994305
$ws.Range("AP9").Value = 15.324218930563593
$ws.Range("AP10").Value = 15.656736395367849
$ws.Range("AP11").Value = 9.4022613603242409
$ws.Range("AP12").Value = 0
$ws.Range("AP13").Value = 11.473445449940511
$ws.Range("AP14").Value = 0
$ws.Range("AP15").Value = 0
$ws.Range("AP16").Value = 6.4035590373244649
$ws.Range("AP17").Value = 0
$ws.Range("AP18").Value = 0

# Match the post-edit selection position recorded in the workbook
$ws.Range("AS7").Select()
